# Apply the scheduled-runner profit recalculations across the Kujata_Profits sheets.
# Each Leve row's market-board columns (H..N) are refreshed from the latest price pull;
# some rows gain/lose the profit cells (M/N) depending on whether HQ/NQ pricing applies.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 10100
$ws.Range("I6").Value = 10100
$ws.Range("K6").Value = 30300
$ws.Range("M6").Value = -30188
$ws.Range("H40").Value = 1361.7273
$ws.Range("I40").Value = 866.6667
$ws.Range("J40").Value = 1547.375
$ws.Range("K40").Value = 866.6667
$ws.Range("L40").Value = 1547.375
$ws.Range("M40").Value = -691.6667
$ws.Range("N40").Value = -1897.375
$ws.Range("H70").Value = 2129.5293
$ws.Range("I70").Value = 2119.2307
$ws.Range("J70").Value = 2163
$ws.Range("K70").Value = 6357.6921
$ws.Range("L70").Value = 6489
$ws.Range("M70").Value = -6087.6921
$ws.Range("N70").Value = -7029
$ws.Range("H73").Value = 2129.5293
$ws.Range("I73").Value = 2119.2307
$ws.Range("J73").Value = 2163
$ws.Range("K73").Value = 6357.6921
$ws.Range("L73").Value = 6489
$ws.Range("M73").Value = -5421.6921
$ws.Range("N73").Value = -8361
$ws.Range("H141").Value = 1519.3684
$ws.Range("I141").Value = 1548.2222
$ws.Range("J141").Value = 1000
$ws.Range("K141").Value = 4644.6666
$ws.Range("L141").Value = 3000
$ws.Range("M141").Value = 535.3334000000004
$ws.Range("N141").Value = -13360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4191.2104
$ws.Range("I32").Value = 4257.9443
$ws.Range("K32").Value = 4257.9443
$ws.Range("M32").Value = -3970.9443
$ws.Range("H74").Value = 1065.6154
$ws.Range("I74").Value = 885.3333
$ws.Range("J74").Value = 1822.8
$ws.Range("K74").Value = 885.3333
$ws.Range("L74").Value = 1822.8
$ws.Range("M74").Value = -11.33330000000001
$ws.Range("N74").Value = -3570.8
$ws.Range("H77").Value = 1065.6154
$ws.Range("I77").Value = 885.3333
$ws.Range("J77").Value = 1822.8
$ws.Range("K77").Value = 4426.6665
$ws.Range("L77").Value = 9114
$ws.Range("M77").Value = -58.66650000000027
$ws.Range("N77").Value = -17850

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 13889722
$ws.Range("I94").Value = 19231392
$ws.Range("K94").Value = 19231392
$ws.Range("M94").Value = -19230941
$ws.Range("H105").Value = 166670030
$ws.Range("I105").Value = 200003280
$ws.Range("J105").Value = 3800
$ws.Range("K105").Value = 200003280
$ws.Range("L105").Value = 3800
$ws.Range("M105").Value = -200001533
$ws.Range("N105").Value = -7294
$ws.Range("H132").Value = 49748
$ws.Range("J132").Value = 49748
$ws.Range("L132").Value = 49748
$ws.Range("N132").Value = -59868

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 12000
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H27").Value = 12000
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 9672.666999999999
$ws.Range("I110").Value = 1527
$ws.Range("J110").Value = 12000
$ws.Range("K110").Value = 4581
$ws.Range("L110").Value = 36000
$ws.Range("M110").Value = -491
$ws.Range("N110").Value = -44180

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("H70").Value = 17312066
$ws.Range("I70").Value = 19234374
$ws.Range("K70").Value = 19234374
$ws.Range("M70").Value = -19234104
$ws.Range("H73").Value = 17312066
$ws.Range("I73").Value = 19234374
$ws.Range("K73").Value = 19234374
$ws.Range("M73").Value = -19233438
$ws.Range("H80").Value = 4814.143
$ws.Range("J80").Value = 6033.1665
$ws.Range("L80").Value = 6033.1665
$ws.Range("N80").Value = -8029.1665
$ws.Range("H83").Value = 4814.143
$ws.Range("J83").Value = 6033.1665
$ws.Range("L83").Value = 30165.8325
$ws.Range("N83").Value = -40149.8325
$ws.Range("H102").Value = 1451.4242
$ws.Range("I102").Value = 1648.1666
$ws.Range("J102").Value = 1215.3334
$ws.Range("K102").Value = 1648.1666
$ws.Range("L102").Value = 1215.3334
$ws.Range("M102").Value = -26.16660000000002
$ws.Range("N102").Value = -4459.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1011.3333
$ws.Range("I22").Value = 601
$ws.Range("J22").Value = 1062.625
$ws.Range("K22").Value = 601
$ws.Range("L22").Value = 1062.625
$ws.Range("M22").Value = -306
$ws.Range("N22").Value = -1652.625
$ws.Range("H25").Value = 527252
$ws.Range("I25").Value = 1002500
$ws.Range("J25").Value = 52004
$ws.Range("K25").Value = 1002500
$ws.Range("L25").Value = 52004
$ws.Range("M25").Value = -1002270
$ws.Range("N25").Value = -52464
$ws.Range("H26").Value = 1200
$ws.Range("I26").Value = 1200
$ws.Range("K26").Value = 1200
$ws.Range("M26").Value = -905
$ws.Range("H27").Value = 1011.3333
$ws.Range("I27").Value = 601
$ws.Range("J27").Value = 1062.625
$ws.Range("K27").Value = 601
$ws.Range("L27").Value = 1062.625
$ws.Range("M27").Value = -494
$ws.Range("N27").Value = -1276.625
$ws.Range("H31").Value = 3736.25
$ws.Range("I31").Value = 2015
$ws.Range("J31").Value = 4310
$ws.Range("K31").Value = 2015
$ws.Range("L31").Value = 4310
$ws.Range("M31").Value = -1767
$ws.Range("N31").Value = -4806
$ws.Range("H132").Value = 50115.145
$ws.Range("I132").Value = 2055.182
$ws.Range("K132").Value = 6165.545999999999
$ws.Range("M132").Value = -3635.545999999999
$ws.Range("H136").Value = 4704.552
$ws.Range("I136").Value = 5973.7
$ws.Range("K136").Value = 17921.1
$ws.Range("M136").Value = -15371.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 1000000000
$ws.Range("I26").Value = 1000000000
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 1000000000
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -999999707
$ws.Range("N26").ClearContents()
$ws.Range("H132").Value = 9409
$ws.Range("I132").Value = 11238.875
$ws.Range("K132").Value = 33716.625
$ws.Range("M132").Value = -31186.625
